$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2..11 correspond to icon numbers 01..10. The cell text in column E
# moves from a GitHub "blob" viewer URL to the raw.githubusercontent.com
# equivalent (same path, minus the "/blob/" segment), and each of those
# cells becomes a live hyperlink to that same URL.
for ($row = 2; $row -le 11; $row++) {
    $n = $row - 1
    $suffix = "{0:D2}" -f $n
    $url = "https://raw.githubusercontent.com/Ing-Aladar-Dukay/CV_Dukay/9561ba175215b42c1ccb1636f483405c3a2d368d/03%20Colors%20icons/color%20" + $suffix + ".png"

    $cell = $ws.Cells.Item($row, 5)
    $cell.Value = $url
    [void]$ws.Hyperlinks.Add($cell, $url)
}

# Selection moves from E11 to A11.
[void]$ws.Range("A11").Select()
